$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 8912.6757053519
$ws.Range("C2").Value = 8292.79958711069
$ws.Range("E2").Value = 4531.51425364618
$ws.Range("F2").Value = 62.3464100315359

$ws.Range("B3").Value = 8531.20075614363
$ws.Range("C3").Value = 7216.50553988944
$ws.Range("E3").Value = 4715.97180717378
$ws.Range("F3").Value = 197.186556127634

$ws.Range("B4").Value = 2957.73248888592
$ws.Range("C4").Value = 4657.6773187927
$ws.Range("E4").Value = 4242.63107761828
$ws.Range("F4").Value = 70.8461831837911

$ws.Range("B5").Value = 2792.158616954
$ws.Range("C5").Value = 4473.17791959754
$ws.Range("E5").Value = 4184.15337453876
$ws.Range("F5").Value = 60.722137255679

$ws.Range("B6").Value = 9067.23413369974
$ws.Range("C6").Value = 7541.9318531657
$ws.Range("E6").Value = 5055.44218596291
$ws.Range("F6").Value = 224.890584963692

$ws.Range("B7").Value = 10173.0852417893
$ws.Range("C7").Value = 8817.74878142309
$ws.Range("E7").Value = 5853.85983809785
$ws.Range("F7").Value = 311.317025813372
